# DVV_original.xlsx update
# Adds two new "kohde" (target) records (200000002C / 200000002D) across the
# four data sheets, and repositions the active sheet/selection to reflect
# where the author ended up after editing ("R9 huon asukk", cell D9/D10).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "R1 rakennus"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 9
$ws1.Range("A9").NumberFormat = "@"
$ws1.Range("A9").Value = "200000002C"
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "398"
$ws1.Range("C9").Value = 20
$ws1.Range("D9").NumberFormat = "@"
$ws1.Range("D9").Value = "39800200030002"
$ws1.Range("F9").NumberFormat = "@"
$ws1.Range("F9").Value = "39800200030002"
$ws1.Range("G9").Value = 1
$ws1.Range("H9").Value = 1
$ws1.Range("J9").NumberFormat = "@"
$ws1.Range("J9").Value = "15230"
$ws1.Range("K9").Value = 19750111
$ws1.Range("L9").Value = 1
$ws1.Range("N9").Value = 80
$ws1.Range("O9").Value = 1
$ws1.Range("Q9").Value = 250
$ws1.Range("S9").NumberFormat = "@"
$ws1.Range("S9").Value = "01"
$ws1.Range("T9").Value = 19780101
$ws1.Range("U9").NumberFormat = "@"
$ws1.Range("U9").Value = "011"
$ws1.Range("V9").Value = 1
$ws1.Range("W9").Value = 1
$ws1.Range("X9").Value = 6765334
$ws1.Range("Y9").Value = 428759

# Row 10
$ws1.Range("A10").NumberFormat = "@"
$ws1.Range("A10").Value = "200000002D"
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "398"
$ws1.Range("C10").Value = 20
$ws1.Range("D10").NumberFormat = "@"
$ws1.Range("D10").Value = "39800200030003"
$ws1.Range("F10").NumberFormat = "@"
$ws1.Range("F10").Value = "39800200030003"
$ws1.Range("G10").Value = 1
$ws1.Range("H10").Value = 1
$ws1.Range("J10").NumberFormat = "@"
$ws1.Range("J10").Value = "15230"
$ws1.Range("K10").Value = 19750111
$ws1.Range("L10").Value = 1
$ws1.Range("N10").Value = 80
$ws1.Range("O10").Value = 1
$ws1.Range("Q10").Value = 250
$ws1.Range("S10").NumberFormat = "@"
$ws1.Range("S10").Value = "01"
$ws1.Range("T10").Value = 19780101
$ws1.Range("U10").NumberFormat = "@"
$ws1.Range("U10").Value = "011"
$ws1.Range("V10").Value = 1
$ws1.Range("W10").Value = 1
$ws1.Range("X10").Value = 6765334
$ws1.Range("Y10").Value = 428759

# ---------------------------------------------------------------------
# Sheet 2: "R3 osoite"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 9
$ws2.Range("A9").NumberFormat = "@"
$ws2.Range("A9").Value = "200000002C"
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "398"
$ws2.Range("C9").Value = 1
$ws2.Range("D9").NumberFormat = "@"
$ws2.Range("D9").Value = "Halmekatu"
$ws2.Range("F9").Value = 14
$ws2.Range("G9").NumberFormat = "@"
$ws2.Range("G9").Value = "15230"
$ws2.Range("H9").NumberFormat = "@"
$ws2.Range("H9").Value = "LAHTI"
$ws2.Range("I9").NumberFormat = "@"
$ws2.Range("I9").Value = "LAHTIS"

# Row 10
$ws2.Range("A10").NumberFormat = "@"
$ws2.Range("A10").Value = "200000002D"
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "398"
$ws2.Range("C10").Value = 1
$ws2.Range("D10").NumberFormat = "@"
$ws2.Range("D10").Value = "Halmekatu"
$ws2.Range("F10").Value = 15
$ws2.Range("G10").NumberFormat = "@"
$ws2.Range("G10").Value = "15230"
$ws2.Range("H10").NumberFormat = "@"
$ws2.Range("H10").Value = "LAHTI"
$ws2.Range("I10").NumberFormat = "@"
$ws2.Range("I10").Value = "LAHTIS"

# ---------------------------------------------------------------------
# Sheet 3: "R4 omistaja"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Row 9
$ws3.Range("A9").NumberFormat = "@"
$ws3.Range("A9").Value = "200000002C"
$ws3.Range("B9").Value = 560
$ws3.Range("C9").NumberFormat = "@"
$ws3.Range("C9").Value = "151046-9873"
$ws3.Range("F9").Value = 20250101
$ws3.Range("G9").Value = 1
$ws3.Range("H9").NumberFormat = "@"
$ws3.Range("H9").Value = "02"
$ws3.Range("I9").NumberFormat = "@"
$ws3.Range("I9").Value = "Lauko Puolikuoma"
$ws3.Range("J9").NumberFormat = "@"
$ws3.Range("J9").Value = "560"
$ws3.Range("K9").NumberFormat = "@"
$ws3.Range("K9").Value = "fi"
$ws3.Range("M9").NumberFormat = "@"
$ws3.Range("M9").Value = "Kirkkoäyrääntie 1d"
$ws3.Range("N9").Value = 16200
$ws3.Range("O9").NumberFormat = "@"
$ws3.Range("O9").Value = "ARTJÄRVI"
$ws3.Range("P9").Value = 19860101

# Row 10
$ws3.Range("A10").NumberFormat = "@"
$ws3.Range("A10").Value = "200000002D"
$ws3.Range("B10").Value = 560
$ws3.Range("C10").NumberFormat = "@"
$ws3.Range("C10").Value = "151046-9874"
$ws3.Range("F10").Value = 20250101
$ws3.Range("G10").Value = 1
$ws3.Range("H10").NumberFormat = "@"
$ws3.Range("H10").Value = "02"
$ws3.Range("I10").NumberFormat = "@"
$ws3.Range("I10").Value = "Kauko Täysikuoma"
$ws3.Range("J10").NumberFormat = "@"
$ws3.Range("J10").Value = "560"
$ws3.Range("K10").NumberFormat = "@"
$ws3.Range("K10").Value = "fi"
$ws3.Range("M10").NumberFormat = "@"
$ws3.Range("M10").Value = "Kirkkoäyrääntie 1d"
$ws3.Range("N10").Value = 16200
$ws3.Range("O10").NumberFormat = "@"
$ws3.Range("O10").Value = "ARTJÄRVI"
$ws3.Range("P10").Value = 19860101

# ---------------------------------------------------------------------
# Sheet 4: "R9 huon asukk"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Row 10
$ws4.Range("A10").NumberFormat = "@"
$ws4.Range("A10").Value = "200000002C"
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = "398"
$ws4.Range("D10").NumberFormat = "@"
$ws4.Range("D10").Value = "001"
$ws4.Range("F10").Value = 1
$ws4.Range("G10").NumberFormat = "@"
$ws4.Range("G10").Value = "130644-0437"
$ws4.Range("H10").NumberFormat = "@"
$ws4.Range("H10").Value = "Marko"
$ws4.Range("I10").NumberFormat = "@"
$ws4.Range("I10").Value = "Poolo"
$ws4.Range("J10").NumberFormat = "@"
$ws4.Range("J10").Value = "Halmekatu 14"
$ws4.Range("K10").NumberFormat = "@"
$ws4.Range("K10").Value = "15230"
$ws4.Range("L10").NumberFormat = "@"
$ws4.Range("L10").Value = "LAHTI"
$ws4.Range("M10").Value = 20010603
$ws4.Range("N10").Value = 2
$ws4.Range("O10").Value = 0

# Row 11
$ws4.Range("A11").NumberFormat = "@"
$ws4.Range("A11").Value = "200000002D"
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "398"
$ws4.Range("D11").NumberFormat = "@"
$ws4.Range("D11").Value = "001"
$ws4.Range("F11").Value = 1
$ws4.Range("G11").NumberFormat = "@"
$ws4.Range("G11").Value = "130694-0534"
$ws4.Range("H11").NumberFormat = "@"
$ws4.Range("H11").Value = "Mahtu"
$ws4.Range("I11").NumberFormat = "@"
$ws4.Range("I11").Value = "Ahtinen"
$ws4.Range("J11").NumberFormat = "@"
$ws4.Range("J11").Value = "Halmekatu 15"
$ws4.Range("K11").NumberFormat = "@"
$ws4.Range("K11").Value = "15230"
$ws4.Range("L11").NumberFormat = "@"
$ws4.Range("L11").Value = "LAHTI"
$ws4.Range("M11").Value = 20010603
$ws4.Range("N11").Value = 1
$ws4.Range("O11").Value = 0

# ---------------------------------------------------------------------
# Selections / active sheet, matching where the author ended up editing.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D10").Select()

$ws2.Activate()
$ws2.Range("I10").Select()

$ws3.Activate()
$ws3.Range("I16").Select()

$ws4.Activate()
$ws4.Range("D9").Select()
